{"js": "// Replace the date header and every two-digit-by-two-digit multiplication\n// prompt in the worksheet with the values from the updated revision.\n// Old text -> new text. Every \"old\" string is unique in the document, so a\n// plain text search/replace (no wildcards) is safe and unambiguous.\nconst pairs = [\n  [\"2025-06-02 Monday\", \"2025-06-03 Tuesday\"],\n  [\"52\u00d789=\", \"25\u00d739=\"],\n  [\"91\u00d796=\", \"44\u00d764=\"],\n  [\"28\u00d789=\", \"59\u00d719=\"],\n  [\"32\u00d748=\", \"87\u00d752=\"],\n  [\"80\u00d721=\", \"31\u00d799=\"],\n  [\"66\u00d755=\", \"28\u00d711=\"],\n  [\"54\u00d750=\", \"83\u00d789=\"],\n  [\"72\u00d726=\", \"12\u00d719=\"],\n  [\"17\u00d782=\", \"28\u00d733=\"],\n  [\"29\u00d732=\", \"44\u00d732=\"],\n  [\"84\u00d741=\", \"79\u00d728=\"],\n  [\"50\u00d767=\", \"83\u00d794=\"],\n  [\"87\u00d712=\", \"37\u00d733=\"],\n  [\"28\u00d771=\", \"19\u00d738=\"],\n  [\"51\u00d723=\", \"43\u00d755=\"],\n  [\"87\u00d762=\", \"33\u00d746=\"],\n  [\"77\u00d771=\", \"55\u00d765=\"],\n  [\"49\u00d717=\", \"65\u00d763=\"],\n  [\"65\u00d767=\", \"63\u00d715=\"],\n  [\"28\u00d731=\", \"61\u00d756=\"],\n  [\"50\u00d743=\", \"18\u00d758=\"],\n  [\"67\u00d757=\", \"83\u00d721=\"],\n  [\"49\u00d754=\", \"67\u00d765=\"],\n  [\"53\u00d734=\", \"25\u00d729=\"],\n  [\"34\u00d775=\", \"46\u00d799=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date header and every two-digit-by-two-digit multiplication\n# prompt in the worksheet with the values from the updated revision.\n# Each \"old\" string occurs exactly once in the document, so a simple\n# Find/Replace (no wildcards) for each pair is unambiguous.\n$pairs = @(\n    @('2025-06-02 Monday', '2025-06-03 Tuesday'),\n    @('52\u00d789=', '25\u00d739='),\n    @('91\u00d796=', '44\u00d764='),\n    @('28\u00d789=', '59\u00d719='),\n    @('32\u00d748=', '87\u00d752='),\n    @('80\u00d721=', '31\u00d799='),\n    @('66\u00d755=', '28\u00d711='),\n    @('54\u00d750=', '83\u00d789='),\n    @('72\u00d726=', '12\u00d719='),\n    @('17\u00d782=', '28\u00d733='),\n    @('29\u00d732=', '44\u00d732='),\n    @('84\u00d741=', '79\u00d728='),\n    @('50\u00d767=', '83\u00d794='),\n    @('87\u00d712=', '37\u00d733='),\n    @('28\u00d771=', '19\u00d738='),\n    @('51\u00d723=', '43\u00d755='),\n    @('87\u00d762=', '33\u00d746='),\n    @('77\u00d771=', '55\u00d765='),\n    @('49\u00d717=', '65\u00d763='),\n    @('65\u00d767=', '63\u00d715='),\n    @('28\u00d731=', '61\u00d756='),\n    @('50\u00d743=', '18\u00d758='),\n    @('67\u00d757=', '83\u00d721='),\n    @('49\u00d754=', '67\u00d765='),\n    @('53\u00d734=', '25\u00d729='),\n    @('34\u00d775=', '46\u00d799=')\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
